$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the two "усач" game-description cells with rephrased versions
# ("усач" -> "один из участников"), which in the authored edit moved the
# two shared strings to the end of the sst and swapped their row order.
# Set B3 first so its shared string is appended before B2's, matching the
# append order recorded in the saved file.
$ws.Range("B3").Value = "один из участников закрывает глаза, девочка рисует на его спине пальцем букву — он должен угадать"
$ws.Range("B2").Value = "игра «повторюха»: один из участников делает движение, девочка повторяет, кто проиграл — читает рэп-куплет из головы"

# Row 2 grew taller to fit the longer wrapped text.
$ws.Rows.Item(2).RowHeight = 62

# Update the active selection as recorded in the saved view state.
$ws.Range("E2").Select()
